$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the ID value in A2 (price values correction)
$ws.Range("A2").Value = 32

# Move/normalize the active selection to A3
$ws.Activate()
$ws.Range("A3").Select()
